# Baseball Catcher: Working on displaying all of the algorithm frame and real frame
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Run 1")
$ws2 = $wb.Worksheets.Item("Run 2")

# --- Sheet2 ("Run 2"): add the algorithm-frame legend next to the real frame ---
$ws2.Range("F1").Value = "catch = 4"
$ws2.Range("F2").Value = "rim = 2"
$ws2.Range("F3").Value = "miss = 0"
$ws2.Range("F4").Value = "uncatchable = x"
$ws2.Range("B6").Value = "  "
$ws2.Columns.Item(6).ColumnWidth = 12.9453125

# --- View / selection updates ---
$ws1.Activate()
$ws1.Range("A1:D31").Select()

$ws2.Activate()
$ws2.Range("A2:D5").Select()
